$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 7-14, keeping only header + rows 2-6
$ws.Rows("7:14").Delete() | Out-Null

# Update remaining data rows with new content
$ws.Range("A2").Value = "LFSN"
$ws.Range("B2").Value = "LFSN1800DIG"
$ws.Range("C2").Value = 46237
$ws.Range("D2").Value = "FIRF GERAES"

$ws.Range("A3").Value = "CDB"
$ws.Range("B3").Value = "CDB725BEF4D"
$ws.Range("C3").Value = 45941
$ws.Range("D3").Value = "BMG SEG"

$ws.Range("A4").Value = "CDB"
$ws.Range("B4").Value = "CDB725BEF4B"
$ws.Range("C4").Value = 45941
$ws.Range("D4").Value = "FIRF GERAES"

$ws.Range("A5").Value = "CDB"
$ws.Range("B5").Value = "CDB725BEF4E"
$ws.Range("C5").Value = 45941
$ws.Range("D5").Value = "FIRF GERAES 30"

$ws.Range("A6").Value = "CDB"
$ws.Range("B6").Value = "CDB725BEF4C"
$ws.Range("C6").Value = 45941
$ws.Range("D6").Value = "HORIZONTE"
